# Refactor: drop the unused "Durée_image_mémo" resource-duration column from
# the question sheet. The G column was only ever populated for one row
# (image-memo, G6=4) and its header (G1); both are removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G1 header text ("Durée_image_mémo") is no longer needed -> clear it but
# keep the cell (and its bold/centered header style) in place.
$ws.Range("G1").ClearContents()

# G6 held the stray duration value (4) for the image-memo row -> remove the
# cell entirely (value + formatting), not just blank its contents.
$ws.Range("G6").Clear()

# Move the active selection to where the author's cursor ended up.
[void]$ws.Range("H13").Select()
